$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$r = $ws.Range("B14")
Write-Output $r.NumberFormat
$r.NumberFormat = "0"
